$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 1229, pushing the existing
# rows 1229-1314 down to 1230-1315.
$ws.Rows.Item(1229).Insert()

$ws.Cells.Item(1229, 1).Value = 10
$ws.Cells.Item(1229, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1229, 3).Value = "La Araucanía"
$ws.Cells.Item(1229, 4).Value = 44610
$ws.Cells.Item(1229, 5).Value = 9
$ws.Cells.Item(1229, 6).Value = 100112020
$ws.Cells.Item(1229, 7).Value = "Tomate"
$ws.Cells.Item(1229, 8).Value = "Larga vida"
$ws.Cells.Item(1229, 9).Value = "Primera"
$ws.Cells.Item(1229, 10).Value = 1500
$ws.Cells.Item(1229, 11).Value = 8000
$ws.Cells.Item(1229, 12).Value = 9000
$ws.Cells.Item(1229, 13).Value = 8533
$ws.Cells.Item(1229, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(1229, 15).Value = "Angol"
$ws.Cells.Item(1229, 16).Value = 474
$ws.Cells.Item(1229, 17).Value = 18
$ws.Cells.Item(1229, 18).Value = "Hortaliza"
